# Refresh the cryptocurrency price/volume figures from the latest scrape.
# "RenderToken" now ranks ahead of "Monero", so every coin from row 35
# ("Monero" before) down to row 51 ("TheGraph" before) shifts down one slot
# (name/link/price/volume all change); "TheGraph" drops off the bottom of
# the list. Rows 2-34 keep the same coin but get refreshed price/volume.
#
# Several "Price" values look numeric (e.g. "587.81") but must stay plain
# text like the rest of column D (thousand-separator values such as
# "63.483.37" can't be numbers anyway, and Excel would otherwise coerce a
# clean numeric-looking string into a real number cell). Set-TextValue
# forces the text number format just long enough to store the literal text,
# then restores the cell's original style so no visible formatting changes.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.483.37'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '2.603.26'
$ws.Range("E3").Value = '  -1.35%  '
Set-TextValue $ws "D5" '587.81'
$ws.Range("E5").Value = '  -3.04%  '
Set-TextValue $ws "D6" '149.24'
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("E9").Value = '  -1.29%  '
Set-TextValue $ws "D10" '5.74'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  -0.80%  '
Set-TextValue $ws "D13" '27.54'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '3.071.23'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").Value = '63.320.41'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("E16").Value = '  +3.03%  '
$ws.Range("D17").Value = '2.640.88'
$ws.Range("E17").Value = '  -0.36%  '
Set-TextValue $ws "D18" '12.06'
$ws.Range("E18").Value = '  -1.30%  '
Set-TextValue $ws "D19" '4.68'
$ws.Range("E19").Value = '  +0.25%  '
Set-TextValue $ws "D20" '343.91'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("E21").Value = '  -2.40%  '
$ws.Range("E22").Value = '  -0.30%  '
Set-TextValue $ws "D23" '66.54'
$ws.Range("E23").Value = '  -0.57%  '
Set-TextValue $ws "D24" '1.72'
Set-TextValue $ws "D25" '9.19'
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").Value = '  -3.88%  '
Set-TextValue $ws "D27" '557.53'
$ws.Range("E27").Value = '  +1.38%  '
Set-TextValue $ws "D28" '8.19'
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  -3.43%  '
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").Value = '0.0₃0848'
$ws.Range("E32").Value = '  -2.43%  '
$ws.Range("E33").Value = '  -1.26%  '
Set-TextValue $ws "D34" '5.27'
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws "D35" '6.08'
$ws.Range("E35").Value = '  -2.01%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws "D36" '165.78'
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws "D37" '0.411'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws "D38" '0.999'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D39" '19.38'
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D40" '1.92'
$ws.Range("E40").Value = '  -5.87%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws "D41" '0.999'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws "D42" '165.82'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D43" '3.98'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws "D44" '22.88'
$ws.Range("E44").Value = '  +5.27%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D45" '0.0578'
$ws.Range("E45").Value = '  -1.98%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D46" '2.09'
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws "D47" '0.631'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D48" '0.0247'
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws "D49" '0.0956'
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D50" '19.03'
$ws.Range("E50").Value = '  -1.99%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0223'
$ws.Range("E51").Value = '  +11.67%  '
